$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.840.44"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.899.17"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "595.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.198"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").Value = "2.897.85"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.426"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.34%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "3.432.26"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").Value = "75.741.58"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000190"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "2.906.15"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +2.99%  "
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "499.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  -5.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "180.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("E42").Value = "  +16.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.343"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.570"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("E51").Value = "  +3.16%  "
